$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new exercise-log record as row 23, copying the formatting
# (date / number / text styles) already used by the preceding rows.
$ws.Range("A22:I22").Copy()
$ws.Range("A23:I23").PasteSpecial(-4122)

$ws.Range("A23").Value = 44008
$ws.Range("B23").Value = 100
$ws.Range("C23").Value = 105
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = "WORSENED"
$ws.Range("F23").Value = 0.1
$ws.Range("G23").Value = 82.6
$ws.Range("H23").Value = 28.9
$ws.Range("I23").Value = "OVERWEIGHT"

# Select the whole sheet, matching the saved state
$ws.Cells.Select()
